$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 10000
$ws.Range("I18").Value = 10000
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 10000
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -9716

$ws.Range("H44").Value = 50000
$ws.Range("J44").Value = 50000
$ws.Range("L44").Value = 50000
$ws.Range("N44").Value = -50924

$ws.Range("H92").Value = 936.4211
$ws.Range("I92").Value = 266.90625
$ws.Range("K92").Value = 266.90625
$ws.Range("M92").Value = 981.09375

$ws.Range("H94").Value = 256.8
$ws.Range("I94").Value = 256.8
$ws.Range("K94").Value = 256.8
$ws.Range("M94").Value = 194.2

$ws.Range("H100").Value = 1971.3334
$ws.Range("J100").Value = 2390.111
$ws.Range("L100").Value = 2390.111
$ws.Range("N100").Value = -3472.111

$ws.Range("H103").Value = 632.9
$ws.Range("I103").Value = 707.1429000000001
$ws.Range("K103").Value = 2121.4287
$ws.Range("M103").Value = -1535.4287

$ws.Range("H106").Value = 166668500
$ws.Range("I106").Value = 200001500
$ws.Range("K106").Value = 200001500
$ws.Range("M106").Value = -200000869

$ws.Range("H125").Value = 9262553
$ws.Range("I125").Value = 2100.875
$ws.Range("J125").Value = 11908396
$ws.Range("K125").Value = 18907.875
$ws.Range("L125").Value = 107175564
$ws.Range("M125").Value = -16447.875
$ws.Range("N125").Value = -107180484

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 428.0909
$ws.Range("I5").Value = 84.75
$ws.Range("J5").Value = 624.2857
$ws.Range("K5").Value = 84.75
$ws.Range("L5").Value = 624.2857
$ws.Range("M5").Value = 27.25
$ws.Range("N5").Value = -848.2857

$ws.Range("H32").Value = 10576.758
$ws.Range("I32").Value = 6018.1953
$ws.Range("K32").Value = 6018.1953
$ws.Range("M32").Value = -5731.1953

$ws.Range("H33").Value = 13166.667
$ws.Range("I33").Value = 13166.667
$ws.Range("K33").Value = 13166.667
$ws.Range("M33").Value = -12837.667

$ws.Range("H61").Value = 3312.7646
$ws.Range("I61").Value = 3146.5518
$ws.Range("K61").Value = 3146.5518
$ws.Range("M61").Value = -2934.5518

$ws.Range("H136").Value = 3312.7646
$ws.Range("I136").Value = 3146.5518
$ws.Range("K136").Value = 9439.6554
$ws.Range("M136").Value = -6889.6554

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 428.0909
$ws.Range("I4").Value = 84.75
$ws.Range("J4").Value = 624.2857
$ws.Range("K4").Value = 84.75
$ws.Range("L4").Value = 624.2857
$ws.Range("M4").Value = 30.25
$ws.Range("N4").Value = -854.2857

$ws.Range("H20").Value = 9806133
$ws.Range("I20").Value = 17545634
$ws.Range("J20").Value = 2763.7334
$ws.Range("K20").Value = 17545634
$ws.Range("L20").Value = 2763.7334
$ws.Range("M20").Value = -17545387
$ws.Range("N20").Value = -3257.7334

$ws.Range("H86").Value = 10011969
$ws.Range("I86").Value = 16684618
$ws.Range("J86").Value = 2996
$ws.Range("K86").Value = 16684618
$ws.Range("L86").Value = 2996
$ws.Range("M86").Value = -16683495
$ws.Range("N86").Value = -5242

$ws.Range("H89").Value = 10011969
$ws.Range("I89").Value = 16684618
$ws.Range("J89").Value = 2996
$ws.Range("K89").Value = 83423090
$ws.Range("L89").Value = 14980
$ws.Range("M89").Value = -83417474
$ws.Range("N89").Value = -26212

$ws.Range("H94").Value = 16682540
$ws.Range("I94").Value = 33335650
$ws.Range("J94").Value = 29429.666
$ws.Range("K94").Value = 33335650
$ws.Range("L94").Value = 29429.666
$ws.Range("M94").Value = -33335199
$ws.Range("N94").Value = -30331.666

$ws.Range("H99").Value = 6213743
$ws.Range("I99").Value = 8405441
$ws.Range("K99").Value = 8405441
$ws.Range("M99").Value = -8403943

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 53820
$ws.Range("J37").Value = 53820
$ws.Range("L37").Value = 161460
$ws.Range("N37").Value = -161684

$ws.Range("H68").Value = 1533.1333
$ws.Range("I68").Value = 701
$ws.Range("J68").Value = 1949.2
$ws.Range("K68").Value = 2103
$ws.Range("L68").Value = 5847.6
$ws.Range("M68").Value = -1292
$ws.Range("N68").Value = -7469.6

$ws.Range("H71").Value = 1533.1333
$ws.Range("I71").Value = 701
$ws.Range("J71").Value = 1949.2
$ws.Range("K71").Value = 6309
$ws.Range("L71").Value = 17542.8
$ws.Range("M71").Value = -2253
$ws.Range("N71").Value = -25654.8

$ws.Range("H75").Value = 4791.857
$ws.Range("J75").Value = 5540.5
$ws.Range("L75").Value = 16621.5
$ws.Range("N75").Value = -18617.5

$ws.Range("H78").Value = 4791.857
$ws.Range("J78").Value = 5540.5
$ws.Range("L78").Value = 49864.5
$ws.Range("N78").Value = -59848.5

$ws.Range("H131").Value = 9262465
$ws.Range("I131").Value = 5954152.5
$ws.Range("J131").Value = 10756541
$ws.Range("K131").Value = 17862457.5
$ws.Range("L131").Value = 32269623
$ws.Range("M131").Value = -17857417.5
$ws.Range("N131").Value = -32279703

$ws.Range("H139").Value = 83335310
$ws.Range("I139").Value = 125001016
$ws.Range("J139").Value = 3900
$ws.Range("K139").Value = 375003048
$ws.Range("L139").Value = 11700
$ws.Range("M139").Value = -374997908
$ws.Range("N139").Value = -21980

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 501000
$ws.Range("I14").Value = 501000
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 501000
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -500832

$ws.Range("H27").Value = 14996.667
$ws.Range("I27").Value = 40000
$ws.Range("J27").Value = 9996
$ws.Range("K27").Value = 40000
$ws.Range("L27").Value = 9996
$ws.Range("M27").Value = -39834
$ws.Range("N27").Value = -10328

$ws.Range("H40").Value = 19999.5
$ws.Range("J40").Value = 19999.5
$ws.Range("L40").Value = 19999.5
$ws.Range("N40").Value = -20301.5

$ws.Range("H97").Value = 567531.4399999999
$ws.Range("I97").Value = 850875.4399999999
$ws.Range("J97").Value = 843.4286
$ws.Range("K97").Value = 850875.4399999999
$ws.Range("L97").Value = 843.4286
$ws.Range("M97").Value = -850379.4399999999
$ws.Range("N97").Value = -1835.4286

$ws.Range("H113").Value = 6175831.5
$ws.Range("I113").Value = 11906607
$ws.Range("J113").Value = 4226.6924
$ws.Range("K113").Value = 11906607
$ws.Range("L113").Value = 4226.6924
$ws.Range("M113").Value = -11904437
$ws.Range("N113").Value = -8566.6924

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1771.4286
$ws.Range("J68").Value = 1999.5
$ws.Range("L68").Value = 1999.5
$ws.Range("N68").Value = -3497.5

$ws.Range("H71").Value = 1771.4286
$ws.Range("J71").Value = 1999.5
$ws.Range("L71").Value = 9997.5
$ws.Range("N71").Value = -17485.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 742912.1
$ws.Range("I54").Value = 2525000
$ws.Range("J54").Value = 30077
$ws.Range("K54").Value = 2525000
$ws.Range("L54").Value = 30077
$ws.Range("M54").Value = -2524480
$ws.Range("N54").Value = -31117

$ws.Range("H62").Value = 7720.886
$ws.Range("I62").Value = 2850
$ws.Range("J62").Value = 8207.975
$ws.Range("K62").Value = 2850
$ws.Range("L62").Value = 8207.975
$ws.Range("M62").Value = -2226
$ws.Range("N62").Value = -9455.975

$ws.Range("H65").Value = 7720.886
$ws.Range("I65").Value = 2850
$ws.Range("J65").Value = 8207.975
$ws.Range("K65").Value = 14250
$ws.Range("L65").Value = 41039.875
$ws.Range("M65").Value = -11130
$ws.Range("N65").Value = -47279.875

$ws.Range("H81").Value = 16676167
$ws.Range("J81").Value = 10916.667
$ws.Range("L81").Value = 21833.334
$ws.Range("N81").Value = -23955.334

$ws.Range("H84").Value = 16676167
$ws.Range("J84").Value = 10916.667
$ws.Range("L84").Value = 109166.67
$ws.Range("N84").Value = -119774.67

$ws.Range("H96").Value = 3390.2083
$ws.Range("I96").Value = 3209.3333
$ws.Range("J96").Value = 3932.8333
$ws.Range("K96").Value = 3209.3333
$ws.Range("L96").Value = 3932.8333
$ws.Range("M96").Value = -1836.3333
$ws.Range("N96").Value = -6678.8333

$ws.Range("H100").Value = 1770.6
$ws.Range("I100").Value = 3255
$ws.Range("J100").Value = 781
$ws.Range("K100").Value = 6510
$ws.Range("L100").Value = 1562
$ws.Range("M100").Value = -5969
$ws.Range("N100").Value = -2644

